# Weekly update: insert this week's new records (2 rows) above the
# existing history for row 164/165, shifting the rest of the history
# down by two rows (old row164 -> row166, ..., old row322 -> row324).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before row 164; this pushes the existing
# rows 164-322 down to 166-324 and keeps all their data/formatting intact.
$ws.Rows("164:165").Insert()

# New row 164 ("Primera" quality) - this week's record, based on the
# values that used to be in row 164 with the price updated (+100).
$ws.Range("A164").Value = 11
$ws.Range("B164").Value = "Vega Monumental Concepción"
$ws.Range("C164").Value = "Bíobío"
$ws.Range("D164").Value = 44880
$ws.Range("E164").Value = 8
$ws.Range("F164").Value = 100112009
$ws.Range("G164").Value = "Acelga"
$ws.Range("H164").Value = "Sin especificar"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 200
$ws.Range("K164").Value = 700
$ws.Range("L164").Value = 800
$ws.Range("M164").Value = 750
$ws.Range("N164").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O164").Value = "Región de Ñuble"
$ws.Range("P164").Value = 750
$ws.Range("Q164").Value = 1
$ws.Range("R164").Value = "Hortaliza"

# New row 165 ("Segunda" quality) - this week's record, based on the
# values that used to be in row 165 with the price updated (+100).
$ws.Range("A165").Value = 11
$ws.Range("B165").Value = "Vega Monumental Concepción"
$ws.Range("C165").Value = "Bíobío"
$ws.Range("D165").Value = 44880
$ws.Range("E165").Value = 8
$ws.Range("F165").Value = 100112009
$ws.Range("G165").Value = "Acelga"
$ws.Range("H165").Value = "Sin especificar"
$ws.Range("I165").Value = "Segunda"
$ws.Range("J165").Value = 100
$ws.Range("K165").Value = 600
$ws.Range("L165").Value = 600
$ws.Range("M165").Value = 600
$ws.Range("N165").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O165").Value = "Región de Ñuble"
$ws.Range("P165").Value = 600
$ws.Range("Q165").Value = 1
$ws.Range("R165").Value = "Hortaliza"
